$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1942446043165468
$ws.Range("C2").Value = 0.5827338129496403
$ws.Range("J2").Value = 0.01079136690647482
$ws.Range("P2").Value = 0.1223021582733813
$ws.Range("S2").Value = 0.08992805755395683
$ws.Range("B3").Value = 0.005917159763313609
$ws.Range("C3").Value = 0.02958579881656805
$ws.Range("J3").Value = 0.03550295857988166
$ws.Range("P3").Value = 0.7100591715976331
$ws.Range("S3").Value = 0.2189349112426036
$ws.Range("J4").Value = 0.04878048780487805
$ws.Range("P4").Value = 0.5853658536585366
$ws.Range("S4").Value = 0.3658536585365854
$ws.Range("B6").Value = 0.04680851063829787
$ws.Range("E6").Value = 0.00425531914893617
$ws.Range("F6").Value = 0.06808510638297872
$ws.Range("J6").Value = 0.2297872340425532
$ws.Range("O6").Value = 0.03404255319148936
$ws.Range("Q6").Value = 0.1659574468085106
$ws.Range("R6").Value = 0.09361702127659574
$ws.Range("S6").Value = 0.3574468085106383
$ws.Range("B7").Value = 0.08808290155440414
$ws.Range("D7").Value = 0.02590673575129534
$ws.Range("F7").Value = 0.05181347150259067
$ws.Range("J7").Value = 0.1450777202072539
$ws.Range("O7").Value = 0.02072538860103627
$ws.Range("Q7").Value = 0.1658031088082902
$ws.Range("R7").Value = 0.06735751295336788
$ws.Range("S7").Value = 0.4352331606217616
$ws.Range("B8").Value = 0.08744394618834081
$ws.Range("D8").Value = 0.02242152466367713
$ws.Range("F8").Value = 0.03139013452914798
$ws.Range("J8").Value = 0.1278026905829596
$ws.Range("O8").Value = 0.02914798206278027
$ws.Range("Q8").Value = 0.1928251121076233
$ws.Range("R8").Value = 0.09641255605381166
$ws.Range("S8").Value = 0.4125560538116592
$ws.Range("B9").Value = 0.08936170212765958
$ws.Range("D9").Value = 0.02553191489361702
$ws.Range("F9").Value = 0.07234042553191489
$ws.Range("J9").Value = 0.1531914893617021
$ws.Range("O9").Value = 0.01276595744680851
$ws.Range("Q9").Value = 0.2553191489361702
$ws.Range("R9").Value = 0.06808510638297872
$ws.Range("S9").Value = 0.323404255319149
$ws.Range("B10").Value = 0.09454545454545454
$ws.Range("D10").Value = 0.01527272727272727
$ws.Range("E10").Value = 0.001454545454545454
$ws.Range("F10").Value = 0.08072727272727273
$ws.Range("J10").Value = 0.136
$ws.Range("O10").Value = 0.01818181818181818
$ws.Range("Q10").Value = 0.1985454545454546
$ws.Range("R10").Value = 0.07563636363636364
$ws.Range("S10").Value = 0.3796363636363637
$ws.Range("G11").Value = 0.1372549019607843
$ws.Range("J11").Value = 0.107843137254902
$ws.Range("K11").Value = 0.2091503267973856
$ws.Range("L11").Value = 0.5196078431372549
$ws.Range("S11").Value = 0.0261437908496732
$ws.Range("G12").Value = 0.6936416184971098
$ws.Range("J12").Value = 0.1907514450867052
$ws.Range("K12").Value = 0.02312138728323699
$ws.Range("L12").Value = 0.05202312138728324
$ws.Range("S12").Value = 0.04046242774566474
$ws.Range("G13").Value = 0.7446808510638298
$ws.Range("J13").Value = 0.2127659574468085
$ws.Range("S13").Value = 0.0425531914893617
$ws.Range("F15").Value = 0.0158102766798419
$ws.Range("H15").Value = 0.1146245059288538
$ws.Range("I15").Value = 0.08300395256916997
$ws.Range("J15").Value = 0.3715415019762846
$ws.Range("K15").Value = 0.05928853754940711
$ws.Range("M15").Value = 0.01185770750988142
$ws.Range("N15").Value = 0.003952569169960474
$ws.Range("O15").Value = 0.07114624505928854
$ws.Range("S15").Value = 0.2687747035573123
$ws.Range("F16").Value = 0.01785714285714286
$ws.Range("H16").Value = 0.1726190476190476
$ws.Range("I16").Value = 0.05952380952380952
$ws.Range("J16").Value = 0.4583333333333333
$ws.Range("K16").Value = 0.1071428571428571
$ws.Range("M16").Value = 0.04166666666666666
$ws.Range("O16").Value = 0.06547619047619048
$ws.Range("S16").Value = 0.07738095238095238
$ws.Range("F17").Value = 0.01033057851239669
$ws.Range("H17").Value = 0.2004132231404959
$ws.Range("I17").Value = 0.115702479338843
$ws.Range("J17").Value = 0.4049586776859504
$ws.Range("K17").Value = 0.09710743801652892
$ws.Range("M17").Value = 0.01239669421487603
$ws.Range("N17").Value = 0.002066115702479339
$ws.Range("O17").Value = 0.05578512396694215
$ws.Range("S17").Value = 0.1012396694214876
$ws.Range("F18").Value = 0.01538461538461539
$ws.Range("H18").Value = 0.1435897435897436
$ws.Range("I18").Value = 0.09230769230769231
$ws.Range("J18").Value = 0.4461538461538462
$ws.Range("K18").Value = 0.1128205128205128
$ws.Range("M18").Value = 0.01025641025641026
$ws.Range("O18").Value = 0.08717948717948718
$ws.Range("S18").Value = 0.09230769230769231
$ws.Range("F19").Value = 0.01808590806330068
$ws.Range("H19").Value = 0.1921627731725697
$ws.Range("I19").Value = 0.1002260738507913
$ws.Range("J19").Value = 0.3760361718161266
$ws.Range("K19").Value = 0.1009796533534288
$ws.Range("M19").Value = 0.02260738507912585
$ws.Range("N19").Value = 0.001507159005275057
$ws.Range("O19").Value = 0.07912584777694047
$ws.Range("S19").Value = 0.1092690278824416

Write-Host "Updated 110 cells"
